$wb = $excel.ActiveWorkbook

# --- survey sheet ---------------------------------------------------------
$survey = $wb.Worksheets.Item("survey")

# Row 10 used to be a half-filled placeholder row (text/select/Table).
# Turn it into the real "select_one yesno" question...
$survey.Range("A10").Value = "select_one yesno"
$survey.Range("B10").Value = "yesno"
$survey.Range("C10").Value = "Yes No"

# ...and add a new row 11 for the "select_one sino" question.
$survey.Rows.Item(11).Insert()
$survey.Range("A11").Value = "select_one sino"
$survey.Range("B11").Value = "sino"
$survey.Range("C11").Value = "Si No"

# --- choices sheet ---------------------------------------------------------
$choices = $wb.Worksheets.Item("choices")

# Add the matching "sino" choice list (same Si/No options as "yesno").
$choices.Rows.Item(14).Insert()
$choices.Rows.Item(14).Insert()
$choices.Range("A14").Value = "sino"
$choices.Range("B14").Value = 0
$choices.Range("C14").Value = "Si"
$choices.Range("A15").Value = "sino"
$choices.Range("B15").Value = 1
$choices.Range("C15").Value = "No"

# --- make "survey" the active tab ------------------------------------------
$survey.Activate()
